$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 21.137664
$ws.Range("H2").Value = 63.412992
$ws.Range("I2").Value = 0.01636081711441431
$ws.Range("J2").Value = 0.01636081711441431
$ws.Range("M2").Value = 35.585194
$ws.Range("N2").Value = 106.755582
$ws.Range("O2").Value = 0.9972091466993565
$ws.Range("P2").Value = 0.9972091466993567
$ws.Range("Q2").Value = 752.1878741468161
$ws.Range("R2").Value = 6769.690867321345
$ws.Range("S2").Value = 0.01631515647396932
$ws.Range("T2").Value = 0.01631515647396932
$ws.Range("G3").Value = 21.137664
$ws.Range("H3").Value = 63.412992
$ws.Range("I3").Value = 0.01636081711441431
$ws.Range("J3").Value = 0.01636081711441431
$ws.Range("M3").Value = 0.093901
$ws.Range("N3").Value = 0.281703
$ws.Range("O3").Value = 0.002631401590341653
$ws.Range("P3").Value = 0.002631401590341654
$ws.Range("Q3").Value = 1.984847787264
$ws.Range("R3").Value = 17.863630085376
$ws.Range("S3").Value = 0.00004305188017415875
$ws.Range("T3").Value = 0.00004305188017415875
$ws.Range("G4").Value = 21.137664
$ws.Range("H4").Value = 63.412992
$ws.Range("I4").Value = 0.01636081711441431
$ws.Range("J4").Value = 0.01636081711441431
$ws.Range("K4").Value = 1.0
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.00569
$ws.Range("N4").Value = 0.01707
$ws.Range("O4").Value = 0.0001594517103017434
$ws.Range("P4").Value = 0.0001594517103017434
$ws.Range("Q4").Value = 0.12027330816
$ws.Range("R4").Value = 1.08245977344
$ws.Range("S4").Value = 0.000002608760270827396
$ws.Range("T4").Value = 0.000002608760270827396
$ws.Range("I5").Value = 0.8979526429041496
$ws.Range("J5").Value = 0.8979526429041496
$ws.Range("M5").Value = 35.585194
$ws.Range("N5").Value = 106.755582
$ws.Range("O5").Value = 0.9972091466993565
$ws.Range("P5").Value = 0.9972091466993567
$ws.Range("Q5").Value = 41283.33473977388
$ws.Range("R5").Value = 371550.012657965
$ws.Range("S5").Value = 0.8954465888068789
$ws.Range("T5").Value = 0.8954465888068791
$ws.Range("I6").Value = 0.8979526429041496
$ws.Range("J6").Value = 0.8979526429041496
$ws.Range("M6").Value = 0.093901
$ws.Range("N6").Value = 0.281703
$ws.Range("O6").Value = 0.002631401590341653
$ws.Range("P6").Value = 0.002631401590341654
$ws.Range("Q6").Value = 108.9370600424296
$ws.Range("R6").Value = 980.4335403818668
$ws.Range("S6").Value = 0.00236287401258947
$ws.Range("T6").Value = 0.00236287401258947
$ws.Range("I7").Value = 0.8979526429041496
$ws.Range("J7").Value = 0.8979526429041496
$ws.Range("K7").Value = 1.0
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.00569
$ws.Range("N7").Value = 0.01707
$ws.Range("O7").Value = 0.0001594517103017434
$ws.Range("P7").Value = 0.0001594517103017434
$ws.Range("Q7").Value = 6.601121091803332
$ws.Range("R7").Value = 59.41008982622999
$ws.Range("S7").Value = 0.0001431800846810373
$ws.Range("T7").Value = 0.0001431800846810373
$ws.Range("G8").Value = 109.7535913333333
$ws.Range("H8").Value = 329.260774
$ws.Range("I8").Value = 0.08495065658413503
$ws.Range("J8").Value = 0.08495065658413503
$ws.Range("M8").Value = 35.585194
$ws.Range("N8").Value = 106.755582
$ws.Range("O8").Value = 0.9972091466993565
$ws.Range("P8").Value = 0.9972091466993567
$ws.Range("Q8").Value = 3905.602839793386
$ws.Range("R8").Value = 35150.42555814047
$ws.Range("S8").Value = 0.08471357176381536
$ws.Range("T8").Value = 0.08471357176381539
$ws.Range("G9").Value = 109.7535913333333
$ws.Range("H9").Value = 329.260774
$ws.Range("I9").Value = 0.08495065658413503
$ws.Range("J9").Value = 0.08495065658413503
$ws.Range("M9").Value = 0.093901
$ws.Range("N9").Value = 0.281703
$ws.Range("O9").Value = 0.002631401590341653
$ws.Range("P9").Value = 0.002631401590341654
$ws.Range("Q9").Value = 10.30597197979133
$ws.Range("R9").Value = 92.75374781812201
$ws.Range("S9").Value = 0.0002235392928360606
$ws.Range("T9").Value = 0.0002235392928360606
$ws.Range("G10").Value = 109.7535913333333
$ws.Range("H10").Value = 329.260774
$ws.Range("I10").Value = 0.08495065658413503
$ws.Range("J10").Value = 0.08495065658413503
$ws.Range("K10").Value = 1.0
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.00569
$ws.Range("N10").Value = 0.01707
$ws.Range("O10").Value = 0.0001594517103017434
$ws.Range("P10").Value = 0.0001594517103017434
$ws.Range("Q10").Value = 0.6244979346866667
$ws.Range("R10").Value = 5.62048141218
$ws.Range("S10").Value = 0.00001354552748359639
$ws.Range("T10").Value = 0.00001354552748359639
$ws.Range("G11").Value = 0.9507383333333334
$ws.Range("H11").Value = 2.852215
$ws.Range("I11").Value = 0.0007358833973011272
$ws.Range("J11").Value = 0.0007358833973011273
$ws.Range("M11").Value = 35.585194
$ws.Range("N11").Value = 106.755582
$ws.Range("O11").Value = 0.9972091466993565
$ws.Range("P11").Value = 0.9972091466993567
$ws.Range("Q11").Value = 33.83220803490334
$ws.Range("R11").Value = 304.48987231413
$ws.Range("S11").Value = 0.0007338296546928806
$ws.Range("T11").Value = 0.0007338296546928808
$ws.Range("G12").Value = 0.9507383333333334
$ws.Range("H12").Value = 2.852215
$ws.Range("I12").Value = 0.0007358833973011272
$ws.Range("J12").Value = 0.0007358833973011273
$ws.Range("M12").Value = 0.093901
$ws.Range("N12").Value = 0.281703
$ws.Range("O12").Value = 0.002631401590341653
$ws.Range("P12").Value = 0.002631401590341654
$ws.Range("Q12").Value = 0.08927528023833334
$ws.Range("R12").Value = 0.803477522145
$ws.Range("S12").Value = 0.000001936404741964205
$ws.Range("T12").Value = 0.000001936404741964205
$ws.Range("G13").Value = 0.9507383333333334
$ws.Range("H13").Value = 2.852215
$ws.Range("I13").Value = 0.0007358833973011272
$ws.Range("J13").Value = 0.0007358833973011273
$ws.Range("K13").Value = 1.0
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.00569
$ws.Range("N13").Value = 0.01707
$ws.Range("O13").Value = 0.0001594517103017434
$ws.Range("P13").Value = 0.0001594517103017434
$ws.Range("Q13").Value = 0.005409701116666666
$ws.Range("R13").Value = 0.04868731004999999
$ws.Range("S13").Value = 0.0000001173378662823221
$ws.Range("T13").Value = 0.0000001173378662823221
